$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Data for columns I (I0) and J (IF) per row, rows 2-13
$data = @{
    2  = @(1, 4)
    3  = @(1, 6)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 6)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 3)
    13 = @(4, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
